$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-15"

# Update the row label for May
$ws.Range("A6").Value = "May (through 05-15)"

# Update May row (row 6) values
$ws.Range("C6").Value = 23
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 21
$ws.Range("F6").Value = 18
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 57

# Update Total row (row 7) values
$ws.Range("C7").Value = 185
$ws.Range("D7").Value = 283
$ws.Range("E7").Value = 267
$ws.Range("F7").Value = 173
$ws.Range("G7").Value = 292
$ws.Range("H7").Value = 578
$ws.Range("I7").Value = 609
